$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (A2:R2)
$row2 = @(5.4, 1, 1, 2, 3, 10, 10, 10, 10, 0.04, 0.03, 0.02, 0.07, 200, 300, 600, 0.2822010582010584, 0.7500000000000001)
# Row 3 values (A3:R3)
$row3 = @(5.6, 1, 1, 3, 3, 10, 10, 10, 10, 0.04, 0.03, 0.02, 0.07, 200, 300, 600, 0.18195978835978824, 0.7567567567567569)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

$ws.Range("A3:R3").Select()
